$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$ws = $wb.Worksheets.Item("P_valores")

$ws.Range("C2").Value = 0.5759828773678666
$ws.Range("D2").Value = 0.6047874795115229
$ws.Range("E2").Value = 0.8865472824986775
$ws.Range("F2").Value = 0.80017960914057

$ws.Range("B3").Value = 0.5759828773678666
$ws.Range("D3").Value = 0.2891591898825707
$ws.Range("E3").Value = 0.4985986387039578
$ws.Range("F3").Value = 0.243650175886734

$ws.Range("B4").Value = 0.6047874795115229
$ws.Range("C4").Value = 0.2891591898825707
$ws.Range("E4").Value = 0.3891791632672612
$ws.Range("F4").Value = 0.6222907548544758

$ws.Range("B5").Value = 0.8865472824986775
$ws.Range("C5").Value = 0.4985986387039578
$ws.Range("D5").Value = 0.3891791632672612
$ws.Range("F5").Value = 0.6201803543818003

$ws.Range("B6").Value = 0.80017960914057
$ws.Range("C6").Value = 0.243650175886734
$ws.Range("D6").Value = 0.6222907548544758
$ws.Range("E6").Value = 0.6201803543818003

# --- Sheet: Estadisticos_DM ---
$ws2 = $wb.Worksheets.Item("Estadisticos_DM")

$ws2.Range("C2").Value = 0.5726206914965135
$ws2.Range("D2").Value = -0.5294506749158775
$ws2.Range("E2").Value = 0.1452973112544494
$ws2.Range("F2").Value = -0.2579751896406907

$ws2.Range("B3").Value = -0.5726206914965135
$ws2.Range("D3").Value = -1.101731997646872
$ws2.Range("E3").Value = -0.6947206117884998
$ws2.Range("F3").Value = -1.21719574849004

$ws2.Range("B4").Value = 0.5294506749158775
$ws2.Range("C4").Value = 1.101731997646872
$ws2.Range("E4").Value = 0.8886993011483927
$ws2.Range("F4").Value = 0.5037206941243789

$ws2.Range("B5").Value = -0.1452973112544494
$ws2.Range("C5").Value = 0.6947206117884998
$ws2.Range("D5").Value = -0.8886993011483927
$ws2.Range("F5").Value = -0.5068041222154503

$ws2.Range("B6").Value = 0.2579751896406907
$ws2.Range("C6").Value = 1.21719574849004
$ws2.Range("D6").Value = -0.5037206941243789
$ws2.Range("E6").Value = 0.5068041222154503
